# MAST-U_home.xlsx — "Add files via upload"
#
# Updates a batch of computed/pasted values on the "transistions" sheet
# (columns H/J and N/P, rows 2-8), turns the two formula cells H8/J8 into
# plain pasted values, and moves the active selection to P14.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("transistions")

# H/J columns (rows 2-8) — mirrored pair of values on each row
$ws.Range("H2").Value = 1.1246
$ws.Range("J2").Value = 1.1246

$ws.Range("H3").Value = 0.2123
$ws.Range("J3").Value = 0.2123

$ws.Range("H4").Value = 0.3896
$ws.Range("J4").Value = 0.3896

$ws.Range("H5").Value = 0.78
$ws.Range("J5").Value = 0.78

$ws.Range("H6").Value = 1.6928
$ws.Range("J6").Value = 1.6928

$ws.Range("H7").Value = 1.2355
$ws.Range("J7").Value = 1.2355

# H8/J8 previously held the formula =0.0518 — replace with a plain value
$ws.Range("H8").Value = 3.0659
$ws.Range("J8").Value = 3.0659

# N/P columns (rows 3-6, 8) — mirrored pair of values on each row
$ws.Range("N3").Value = 1.1808
$ws.Range("P3").Value = 1.1808

$ws.Range("N4").Value = 0.2224
$ws.Range("P4").Value = 0.2224

$ws.Range("N5").Value = 0.6658
$ws.Range("P5").Value = 0.6658

$ws.Range("N6").Value = 1.0856
$ws.Range("P6").Value = 1.0856

$ws.Range("N8").Value = 0.3203
$ws.Range("P8").Value = 0.3203

# Move the active selection on the "transistions" sheet to P14
$ws.Range("P14").Select()
